$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 = "president", formatted like the other header cells (bold/centered via E1's style)
$f1 = $ws.Range("F1")
$f1.Value = "president"
$ws.Range("E1").Copy()
$f1.PasteSpecial(-4122)  # xlPasteFormats

# Data cells F2:F229 = "Clinton" for every existing data row
$ws.Range("F2:F229").Value = "Clinton"

Write-Output "done"
